# Adjusted to show kick-off times and not broadcast times
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update event times (kick-off time instead of broadcast time) ---
$ws.Range("F3").Value  = 0.67708333333333337   # Horse Racing: Grand National       13:00 -> 16:15
$ws.Range("F7").Value  = 0.58333333333333337   # Tennis: Wimbledon Men's Final      13:00 -> 14:00
$ws.Range("F8").Value  = 0.58333333333333337   # Tennis: Wimbledon Women's Final    13:30 -> 14:00
$ws.Range("F10").Value = 0.71875               # Soccer: FA Cup Final               15:15 -> 17:15
$ws.Range("F13").Value = 0.66666666666666663   # Soccer: League Cup Final           15:00 -> 16:00
$ws.Range("F21").Value = 0.60416666666666663   # Rugby: England vs. New Zealand     16:30 -> 14:30

# --- Update the view / selection so the window is scrolled back to the
#     top of the sheet and the previously-selected cell follows the
#     two-row shift caused by the edits above (F21 -> F19). ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.Left = 60795
$excel.ActiveWindow.Top = 1260

$ws.Range("F19").Select()
